$d = $word.ActiveDocument

# Update the date line (first paragraph, outside the table)
$d.Content.Find.Execute("2025-04-27 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-28 Monday", 2)

# Update each answer cell in the practice table (row-major order, matches
# the reading order of the document; one pair of values, "72-20=52", is
# not unique across the table so Find/Replace alone can't disambiguate --
# address cells directly by (row, col) instead).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "29+62=91"
$t.Cell(1, 2).Range.Text = "75-41=34"
$t.Cell(1, 3).Range.Text = "14+45=59"
$t.Cell(1, 4).Range.Text = "32+50=82"
$t.Cell(1, 5).Range.Text = "25-9=16"
$t.Cell(2, 1).Range.Text = "48+24=72"
$t.Cell(2, 2).Range.Text = "70-51=19"
$t.Cell(2, 3).Range.Text = "31-10=21"
$t.Cell(2, 4).Range.Text = "61+35=96"
$t.Cell(2, 5).Range.Text = "35-35=0"
$t.Cell(3, 1).Range.Text = "62-39=23"
$t.Cell(3, 2).Range.Text = "22-12=10"
$t.Cell(3, 3).Range.Text = "32+18=50"
$t.Cell(3, 4).Range.Text = "17+58=75"
$t.Cell(3, 5).Range.Text = "88-80=8"
$t.Cell(4, 1).Range.Text = "91-5=86"
$t.Cell(4, 2).Range.Text = "76-51=25"
$t.Cell(4, 3).Range.Text = "50+10=60"
$t.Cell(4, 4).Range.Text = "13+61=74"
$t.Cell(4, 5).Range.Text = "60+15=75"
$t.Cell(5, 1).Range.Text = "41-10=31"
$t.Cell(5, 2).Range.Text = "0+70=70"
$t.Cell(5, 3).Range.Text = "85+0=85"
$t.Cell(5, 4).Range.Text = "55+3=58"
$t.Cell(5, 5).Range.Text = "18+55=73"
$t.Cell(6, 1).Range.Text = "42+18=60"
$t.Cell(6, 2).Range.Text = "7+71=78"
$t.Cell(6, 3).Range.Text = "0+78=78"
$t.Cell(6, 4).Range.Text = "37-13=24"
$t.Cell(6, 5).Range.Text = "90-1=89"
$t.Cell(7, 1).Range.Text = "13+84=97"
$t.Cell(7, 2).Range.Text = "71-34=37"
$t.Cell(7, 3).Range.Text = "70+3=73"
$t.Cell(7, 4).Range.Text = "18+1=19"
$t.Cell(7, 5).Range.Text = "97-78=19"
$t.Cell(8, 1).Range.Text = "10+57=67"
$t.Cell(8, 2).Range.Text = "3+21=24"
$t.Cell(8, 3).Range.Text = "25+60=85"
$t.Cell(8, 4).Range.Text = "56-30=26"
$t.Cell(8, 5).Range.Text = "12+61=73"
$t.Cell(9, 1).Range.Text = "95-80=15"
$t.Cell(9, 2).Range.Text = "54-39=15"
$t.Cell(9, 3).Range.Text = "64-32=32"
$t.Cell(9, 4).Range.Text = "90-46=44"
$t.Cell(9, 5).Range.Text = "85-44=41"
$t.Cell(10, 1).Range.Text = "31+26=57"
$t.Cell(10, 2).Range.Text = "92-83=9"
$t.Cell(10, 3).Range.Text = "48+20=68"
$t.Cell(10, 4).Range.Text = "79-6=73"
$t.Cell(10, 5).Range.Text = "89-71=18"
$t.Cell(11, 1).Range.Text = "12-7=5"
$t.Cell(11, 2).Range.Text = "4+80=84"
$t.Cell(11, 3).Range.Text = "98-7=91"
$t.Cell(11, 4).Range.Text = "93-92=1"
$t.Cell(11, 5).Range.Text = "44+45=89"
$t.Cell(12, 1).Range.Text = "53+11=64"
$t.Cell(12, 2).Range.Text = "34-28=6"
$t.Cell(12, 3).Range.Text = "48+27=75"
$t.Cell(12, 4).Range.Text = "34-34=0"
$t.Cell(12, 5).Range.Text = "19+24=43"
$t.Cell(13, 1).Range.Text = "78-64=14"
$t.Cell(13, 2).Range.Text = "37+33=70"
$t.Cell(13, 3).Range.Text = "45-8=37"
$t.Cell(13, 4).Range.Text = "89-0=89"
$t.Cell(13, 5).Range.Text = "90-89=1"
$t.Cell(14, 1).Range.Text = "96-85=11"
$t.Cell(14, 2).Range.Text = "76-34=42"
$t.Cell(14, 3).Range.Text = "82-8=74"
$t.Cell(14, 4).Range.Text = "57-53=4"
$t.Cell(14, 5).Range.Text = "83-17=66"
$t.Cell(15, 1).Range.Text = "96-88=8"
$t.Cell(15, 2).Range.Text = "41+41=82"
$t.Cell(15, 3).Range.Text = "52+38=90"
$t.Cell(15, 4).Range.Text = "7+18=25"
$t.Cell(15, 5).Range.Text = "46+22=68"
$t.Cell(16, 1).Range.Text = "74-42=32"
$t.Cell(16, 2).Range.Text = "16+5=21"
$t.Cell(16, 3).Range.Text = "99-13=86"
$t.Cell(16, 4).Range.Text = "64+16=80"
$t.Cell(16, 5).Range.Text = "22+1=23"
$t.Cell(17, 1).Range.Text = "43-12=31"
$t.Cell(17, 2).Range.Text = "7+87=94"
$t.Cell(17, 3).Range.Text = "0+14=14"
$t.Cell(17, 4).Range.Text = "23-19=4"
$t.Cell(17, 5).Range.Text = "60+29=89"
$t.Cell(18, 1).Range.Text = "0+10=10"
$t.Cell(18, 2).Range.Text = "93-83=10"
$t.Cell(18, 3).Range.Text = "68-59=9"
$t.Cell(18, 4).Range.Text = "40+52=92"
$t.Cell(18, 5).Range.Text = "19+62=81"
$t.Cell(19, 1).Range.Text = "84-32=52"
$t.Cell(19, 2).Range.Text = "73-41=32"
$t.Cell(19, 3).Range.Text = "86+0=86"
$t.Cell(19, 4).Range.Text = "12+40=52"
$t.Cell(19, 5).Range.Text = "28+27=55"
$t.Cell(20, 1).Range.Text = "12+18=30"
$t.Cell(20, 2).Range.Text = "94+2=96"
$t.Cell(20, 3).Range.Text = "66+32=98"
$t.Cell(20, 4).Range.Text = "36-7=29"
$t.Cell(20, 5).Range.Text = "45-14=31"